# Tidskalkyle.xlsx update:
#  - add a new "102819_D_Bakside" measurement row to "Raw Data"
#  - add the corresponding aggregated "102819" row to "Results"
#  - leave the "Raw Data" sheet selected/active (instead of "Results")

$wb = $excel.ActiveWorkbook

$wsRaw = $wb.Worksheets.Item("Raw Data")
$wsResults = $wb.Worksheets.Item("Results")

# --- Raw Data: new product-code / time(s) entry in row 7 ---
$wsRaw.Range("A7").Value = "102819_D_Bakside"
$wsRaw.Range("B7").Value = 1213

# --- Results: new aggregated row 5 (product code + formatted total time) ---
# Written via TEXT()+PasteSpecial so the numeric-looking "102819" lands as a
# plain text/shared-string value (matching the rest of the column) instead
# of Excel auto-coercing it to a number.
$wsResults.Range("A5").Formula = "=TEXT(102819,""0"")"
$wsResults.Range("A5").Copy() | Out-Null
$wsResults.Range("A5").PasteSpecial(-4163) | Out-Null
$wsResults.Range("B5").Value = "0:20:13"

# --- Selection / active sheet bookkeeping ---
# Select Results!A7 first (it's currently the active sheet), then switch to
# and select Raw Data!A9, leaving "Raw Data" as the active/selected tab.
$wsResults.Range("A7").Select() | Out-Null
$wsRaw.Range("A9").Select() | Out-Null
